# Append new job listing rows (87-90) to the dice jobs list worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Senior IT Software Developer (Go, AWS, Terraform)", "https://www.dice.com/job-detail/f4110968-dd81-4d78-b8f8-5bf7221dab35", "Highland Heights, Ohio", "Contract", "$50 - $55", "HYR Global Source Inc"),
    @("Golang Developer", "https://www.dice.com/job-detail/3276f88e-105c-43fa-afda-bb67b0d4758c", "West Chester, Pennsylvania", "Full-time, Part-time, Third Party, Contract", "Depends on Experience", "NasTech Global, Inc."),
    @("Golang Developer - Remote", "https://www.dice.com/job-detail/79129310-8edc-410f-aa6e-03ca94028af5", "Remote", "Full-time, Contract", "Depends on Experience", "InfiCare Technologies"),
    @("LeadGoLang Developer", "https://www.dice.com/job-detail/d34d8c56-5dd4-4269-84f0-5517893013e3", "Remote", "Contract, Third Party", "Depends on Experience", "Stellar IT Solution")
)

$startRow = 87
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}
